$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "FTQ2MN"
$ws.Range("B71").Value = "Grasa para fusor"
$ws.Range("C71").Value = "20gr"
$ws.Range("D71").Value = 20000
$ws.Range("E71").Value = 75000
$ws.Range("F71").Value = 4
$ws.Range("G71").Value = 5
$ws.Range("H71").Formula = "=(E71-D71)*G71"
$ws.Range("I71").Formula = "=D71*F71"
$ws.Range("J71").Value = 80000
